$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 2213.75
$ws.Range("I96").Value = 1713
$ws.Range("J96").Value = 2714.5
$ws.Range("K96").Value = 5139
$ws.Range("L96").Value = 8143.5
$ws.Range("M96").Value = -3766
$ws.Range("N96").Value = -10889.5

$ws.Range("H113").Value = 3918.5
$ws.Range("I113").Value = 4053.75
$ws.Range("K113").Value = 4053.75
$ws.Range("M113").Value = -799.75

$ws.Range("H129").Value = 888.16327
$ws.Range("I129").Value = 526.8889
$ws.Range("J129").Value = 924.69666
$ws.Range("K129").Value = 1580.6667
$ws.Range("L129").Value = 2774.08998
$ws.Range("M129").Value = 3419.3333
$ws.Range("N129").Value = -12774.08998

$ws.Range("H135").Value = 28911.861
$ws.Range("I135").Value = 36369.32
$ws.Range("J135").Value = 2810.75
$ws.Range("K135").Value = 327323.88
$ws.Range("L135").Value = 25296.75
$ws.Range("M135").Value = -324788.88
$ws.Range("N135").Value = -30366.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1382.0952
$ws.Range("I2").Value = 1094.4375
$ws.Range("J2").Value = 2302.6
$ws.Range("K2").Value = 1094.4375
$ws.Range("L2").Value = 2302.6
$ws.Range("M2").Value = -981.4375
$ws.Range("N2").Value = -2528.6

$ws.Range("H116").Value = 1382.0952
$ws.Range("I116").Value = 1094.4375
$ws.Range("J116").Value = 2302.6
$ws.Range("K116").Value = 1094.4375
$ws.Range("L116").Value = 2302.6
$ws.Range("M116").Value = 1199.5625
$ws.Range("N116").Value = -6890.6

$ws.Range("H125").Value = 60000
$ws.Range("J125").Value = 60000
$ws.Range("L125").Value = 60000
$ws.Range("N125").Value = -69840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1382.0952
$ws.Range("I3").Value = 1094.4375
$ws.Range("J3").Value = 2302.6
$ws.Range("K3").Value = 1094.4375
$ws.Range("L3").Value = 2302.6
$ws.Range("M3").Value = -980.4375
$ws.Range("N3").Value = -2530.6

$ws.Range("H122").Value = 52105.26
$ws.Range("J122").Value = 52105.26
$ws.Range("L122").Value = 52105.26
$ws.Range("N122").Value = -61905.26

$ws.Range("H125").Value = 30000
$ws.Range("J125").Value = 30000
$ws.Range("L125").Value = 30000
$ws.Range("N125").Value = -39840

$ws.Range("H126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

$ws.Range("H129").Value = 49332.668
$ws.Range("J129").Value = 49332.668
$ws.Range("L129").Value = 49332.668
$ws.Range("N129").Value = -59332.668

$ws.Range("H130").Value = 50780
$ws.Range("J130").Value = 50780
$ws.Range("L130").Value = 50780
$ws.Range("N130").Value = -60820

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 51000
$ws.Range("J20").Value = 51000
$ws.Range("L20").Value = 51000
$ws.Range("N20").Value = -51472

$ws.Range("H30").Value = 51000
$ws.Range("J30").Value = 51000
$ws.Range("L30").Value = 51000
$ws.Range("N30").Value = -51182

$ws.Range("H31").Value = 2646.7144
$ws.Range("I31").Value = 1255.8928
$ws.Range("J31").Value = 8210
$ws.Range("K31").Value = 1255.8928
$ws.Range("L31").Value = 8210
$ws.Range("M31").Value = -960.8928000000001
$ws.Range("N31").Value = -8800

$ws.Range("H34").Value = 2646.7144
$ws.Range("I34").Value = 1255.8928
$ws.Range("J34").Value = 8210
$ws.Range("K34").Value = 1255.8928
$ws.Range("L34").Value = 8210
$ws.Range("M34").Value = -1053.8928
$ws.Range("N34").Value = -8614

$ws.Range("H122").Value = 1708.0416
$ws.Range("I122").Value = 1378.579
$ws.Range("J122").Value = 2960
$ws.Range("K122").Value = 4135.737
$ws.Range("L122").Value = 8880
$ws.Range("M122").Value = -1685.737
$ws.Range("N122").Value = -13780

$ws.Range("H127").Value = 32500
$ws.Range("J127").Value = 32500
$ws.Range("L127").Value = 32500
$ws.Range("N127").Value = -42420

$ws.Range("H128").Value = 51000
$ws.Range("J128").Value = 51000
$ws.Range("L128").Value = 51000
$ws.Range("N128").Value = -60960

$ws.Range("H130").Value = 65000
$ws.Range("J130").Value = 65000
$ws.Range("L130").Value = 65000
$ws.Range("N130").Value = -75040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 23299.387
$ws.Range("J5").Value = 996.1667
$ws.Range("L5").Value = 2988.5001
$ws.Range("N5").Value = -3212.5001

$ws.Range("H33").Value = 306.92307
$ws.Range("I33").Value = 265
$ws.Range("K33").Value = 1590
$ws.Range("M33").Value = -1307

$ws.Range("H122").Value = 1182.3043
$ws.Range("I122").Value = 386.2857
$ws.Range("J122").Value = 1530.5625
$ws.Range("K122").Value = 3476.5713
$ws.Range("L122").Value = 13775.0625
$ws.Range("M122").Value = -1026.5713
$ws.Range("N122").Value = -18675.0625

$ws.Range("H132").Value = 1135.9333
$ws.Range("I132").Value = 559.75
$ws.Range("J132").Value = 1345.4546
$ws.Range("K132").Value = 5037.75
$ws.Range("L132").Value = 12109.0914
$ws.Range("M132").Value = -2507.75
$ws.Range("N132").Value = -17169.0914

$ws.Range("H135").Value = 23299.387
$ws.Range("J135").Value = 996.1667
$ws.Range("L135").Value = 8965.5003
$ws.Range("N135").Value = -14035.5003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3717.3333
$ws.Range("I61").Value = 3860.8
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 3860.8
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -3658.8
$ws.Range("N61").Value = -3404

$ws.Range("H74").Value = 3348399
$ws.Range("I74").Value = 10000197
$ws.Range("J74").Value = 22500
$ws.Range("K74").Value = 10000197
$ws.Range("L74").Value = 22500
$ws.Range("M74").Value = -9999199
$ws.Range("N74").Value = -24496

$ws.Range("H77").Value = 3348399
$ws.Range("I77").Value = 10000197
$ws.Range("J77").Value = 22500
$ws.Range("K77").Value = 30000591
$ws.Range("L77").Value = 67500
$ws.Range("M77").Value = -29995599
$ws.Range("N77").Value = -77484

$ws.Range("H113").Value = 3717.3333
$ws.Range("I113").Value = 3860.8
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 3860.8
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -1690.8
$ws.Range("N113").Value = -7340

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 11985
$ws.Range("J49").Value = 11985
$ws.Range("L49").Value = 11985
$ws.Range("N49").Value = -12445

$ws.Range("H75").Value = 22000
$ws.Range("J75").Value = 22000
$ws.Range("L75").Value = 22000
$ws.Range("N75").Value = -23872

$ws.Range("H78").Value = 22000
$ws.Range("J78").Value = 22000
$ws.Range("L78").Value = 66000
$ws.Range("N78").Value = -75360

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
